$d = $word.ActiveDocument

# =========================================================================
# Edit 1: paragraph 1 - "This is a Microsoft word document." gains two
# trailing spaces, then a red parenthetical note is appended. The note is
# built up as three separate Insert operations (matching how the original
# authoring session produced three runs of red text).
# =========================================================================
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "This is a Microsoft word document.  ", 2)

$enDash = [char]0x2013

$r1 = $d.Paragraphs(1).Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Collapse(0)
$r1.InsertAfter("(This is a change " + $enDash + " Ve")
$r1.Font.Color = 255

$r2 = $d.Paragraphs(1).Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Collapse(0)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

$r3 = $d.Paragraphs(1).Range
$r3.MoveEnd(1, -1) | Out-Null
$r3.Collapse(0)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# =========================================================================
# Edit 2: delete the trailing "...ank God almighty, we are free at last."
# paragraph entirely (content + paragraph mark), so the poem's final
# stanza line becomes the last paragraph in the body.
# =========================================================================
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count).Range
$last.Delete()

# =========================================================================
# Edit 3: drop the now-unused custom styles that the reverted/"main
# branch" version of the file never defined (Heading2/Heading4 and their
# linked char styles, plus a handful of web-paste leftovers). None of
# these styles are referenced by any remaining paragraph or run, so
# deleting them is purely cleanup. Styles are removed starting from the
# end of the styles collection and working backwards, which keeps each
# style's index stable for the remaining deletions.
# =========================================================================
$stylesToDelete = @(
    "podcast-toolssubscribe-links",
    "generic-title",
    "subscribe-more-info",
    "subscribe",
    "audio-tool",
    "Heading4Char",
    "Heading2Char",
    "Hyperlink",
    "apple-converted-space",
    "Heading4",
    "Heading2"
)
foreach ($styleName in $stylesToDelete) {
    $style = $d.Styles($styleName)
    $style.Delete()
}
